# The slide number / module label lives in a textbox on the Slide Master
# (it's reused on every slide's footer area), not on an individual slide.
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$shape = $master.Shapes.Item("TextBox 8")
$shape.TextFrame.TextRange.Text = "Module 1 "
